$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new country rows, top-to-bottom, shifting existing rows down.
# Row positions below are expressed in terms of the *current* sheet state at
# the moment each Insert() runs (i.e. after any prior inserts already shifted
# things down), so that the final layout matches the target:
#   4  Egypte / EG          (new, before Erythrée)
#   17 Centafrique / CE     (new, before what is then row 18)
#   18 Côte Ivoire / CI     (new, before Guinée)

$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Egypte"
$ws.Range("B4").Value = "EG"

$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "Centafrique"
$ws.Range("B17").Value = "CE"

$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "Côte Ivoire"
$ws.Range("B18").Value = "CI"
